$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply formatting to the new rows (8-10) by copying the existing
# formats from row 2 (column A style) and row 1 (column B style), so the
# new cells reuse the workbook's existing style indices instead of minting
# new ones. ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A8:A10").PasteSpecial(-4122) | Out-Null

$ws.Range("B1").Copy() | Out-Null
$ws.Range("B8:B10").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Match the row height used by the rest of the sheet.
$ws.Rows("8:10").RowHeight = 15.75

# --- Fill in the new values. Order matches how the shared-string table
# gets built: column A for rows 8-9, then column B for rows 8-9, then the
# whole of row 10. ---
$ws.Range("A8").Value = "list 1"
$ws.Range("A9").Value = "list 2"
$ws.Range("B8").Value = "[a,b,c]"
$ws.Range("B9").Value = '["c","d","e"]'
$ws.Range("A10").Value = "list 3"
$ws.Range("B10").Value = "['f', 'g', 'h']"

$ws.Range("C8").Value = 7
$ws.Range("C9").Value = 8
$ws.Range("C10").Value = 9

# Match the author's final selection.
$ws.Range("B20").Select() | Out-Null
